$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1.46
$ws.Range("F4").Value = 1.08
$ws.Range("G4").Value = 0.91
$ws.Range("D6").Value = 1.54
$ws.Range("G6").Value = 1.05
$ws.Range("D7").Value = 1.77
$ws.Range("F7").Value = 1.45
